$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.357.73"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "3.459.48"
$ws.Range("E3").Value = "  -1.52%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.90"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.04"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.607"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").Value = "3.458.79"
$ws.Range("E9").Value = "  -1.54%  "
$ws.Range("E10").Value = "  -2.20%  "
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.417"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "4.057.45"
$ws.Range("E13").Value = "  -1.54%  "
$ws.Range("E14").Value = "  +1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.59"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.23%  "
$ws.Range("D16").Value = "66.260.59"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "3.452.33"
$ws.Range("E18").Value = "  -1.82%  "
$ws.Range("E19").Value = "  -2.27%  "
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.62"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("E25").Value = "  +0.21%  "
$ws.Range("E26").Value = "  +3.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.32%  "
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -2.40%  "
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.72"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.28"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.59%  "
$ws.Range("B35").Value = "Aptos"
$ws.Range("C35").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.02%  "
$ws.Range("E36").Value = "  -0.31%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "161.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.884"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "28.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("E41").Value = "  +0.64%  "
$ws.Range("D42").Value = "2.776.62"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("E44").Value = "  -2.13%  "
$ws.Range("E45").Value = "  -1.80%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "337.68"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  -2.07%  "
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.96%  "
